# Edit the EBEWE compliance table:
#  - Rename "COMPLIANCE DUE DATE" header -> "INITIAL COMPLIANCE DUE DATE"
#  - Rename "SUBSEQUENT DUE DATE" header -> "INITIAL COMPARATIVE PERIOD" and
#    replace its "Every 5 years thereafter" values with explicit comparative
#    date ranges.
#  - Add a new "NEXT COMPLIANCE DUE DATE" column (D) with the next due dates,
#    formatted as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column D. The sheet already carried a leftover width
# declaration for a (no-longer-present) 4th column; inserting shifts that
# declaration out to column E, where it belongs, and frees up column D for
# the new "NEXT COMPLIANCE DUE DATE" data below.
$ws.Columns.Item(4).Insert()

# --- Column B header ---
$ws.Range("B1").Value = "INITIAL COMPLIANCE DUE DATE"

# --- Column C header + data (was SUBSEQUENT DUE DATE / Every 5 years thereafter) ---
$ws.Range("C1").Value = "INITIAL COMPARATIVE PERIOD"
$ws.Range("C2").Value = "Dec 1, 2016 - Dec 1, 2021"
$ws.Range("C3").Value = "Dec 1, 2017 - Dec 1, 2022"
$ws.Range("C4").Value = "Dec 1, 2018 - Dec 1, 2023"
$ws.Range("C5").Value = "Dec 1, 2019 - Dec 1, 2024"
$ws.Range("C6").Value = "Dec 1, 2020 - Dec 1, 2025"

# --- New column D: NEXT COMPLIANCE DUE DATE ---
# Copy A1's formatting (bold font, borders, centered alignment) onto the new
# header cell, then overwrite its value/content.
$ws.Range("A1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "NEXT COMPLIANCE DUE DATE"

$ws.Range("D2").Value = "Dec 1, 2026"
$ws.Range("D3").Value = "Dec 1, 2027"
$ws.Range("D4").Value = "Dec 1, 2028"
$ws.Range("D5").Value = "Dec 1, 2029"
$ws.Range("D6").Value = "Dec 1, 2030"

$ws.Range("D2:D6").NumberFormat = "@"

# --- Column widths: C and D both end up the same width as column B ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# --- Selection, matching the saved view state ---
$ws.Range("D3").Select() | Out-Null
